$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Filtered Repositories"

# Update header cell B1
$ws.Range("B1").Value = "Full Name"

# Remove the data row (row 2) entirely
$ws.Rows.Item(2).Delete()
